$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.234.85'
$ws.Range("E2").Value = '  -0.48%  '
$ws.Range("D3").Value = '1.829.23'
$ws.Range("E3").Value = '  -0.67%  '
$ws.Range("E4").Value = '  +0.31%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.83'
$ws.Range("E5").Value = '  -1.31%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6025'
$ws.Range("E6").Value = '  -3.85%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.003'
$ws.Range("E7").Value = '  +0.24%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07055'
$ws.Range("E8").Value = '  -5.04%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.2791'
$ws.Range("E9").Value = '  -3.50%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '23.59'
$ws.Range("E10").Value = '  -5.44%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07654'
$ws.Range("E11").Value = '  -0.93%  '
$ws.Range("D12").Value = '1.835.28'
$ws.Range("E12").Value = '  -0.32%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.786'
$ws.Range("E13").Value = '  -3.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6278'
$ws.Range("E14").Value = '  -6.88%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.000009866'
$ws.Range("E15").Value = '  -3.93%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '79.04'
$ws.Range("E16").Value = '  -3.30%  '
$ws.Range("D17").Value = '29.233.29'
$ws.Range("E17").Value = '  -0.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '5.837'
$ws.Range("E18").Value = '  -5.96%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '224.16'
$ws.Range("E19").Value = '  -3.99%  '
$ws.Range("E20").Value = '  +0.21%  '
$ws.Range("E21").Value = '  -5.01%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.006'
$ws.Range("E22").Value = '  -3.96%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.31%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '156.41'
$ws.Range("E24").Value = '  -0.96%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1303'
$ws.Range("E25").Value = '  -3.10%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.975'
$ws.Range("E26").Value = '  -6.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '16.60'
$ws.Range("E27").Value = '  -4.12%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.477'
$ws.Range("E28").Value = '  +0.76%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.06493'
$ws.Range("E29").Value = '  -10.14%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.447'
$ws.Range("E30").Value = '  -2.17%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '3.842'
$ws.Range("E31").Value = '  -4.69%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.792'
$ws.Range("E32").Value = '  -6.13%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.107'
$ws.Range("E33").Value = '  -2.86%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.722'
$ws.Range("E34").Value = '  -5.17%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6465'
$ws.Range("E35").Value = '  -7.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.545'
$ws.Range("D37").Value = '1.214.54'
$ws.Range("E37").Value = '  -1.55%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.738'
$ws.Range("E38").Value = '  -2.72%  '
$ws.Range("E39").Value = '  -5.23%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.559'
$ws.Range("E40").Value = '  -5.35%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.8973'
$ws.Range("E41").Value = '  -6.85%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.003'
$ws.Range("E42").Value = '  +0.25%  '
$ws.Range("D43").Value = '1.992.55'
$ws.Range("E43").Value = '  -0.59%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '100.29'
$ws.Range("E44").Value = '  -0.56%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '62.66'
$ws.Range("E45").Value = '  -4.15%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.00000000116'
$ws.Range("E46").Value = '  -3.23%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '8.568'
$ws.Range("E47").Value = '  -3.53%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.581'
$ws.Range("E48").Value = '  -7.82%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.4549'
$ws.Range("E49").Value = '  -0.55%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.05496'
$ws.Range("E50").Value = '  -2.85%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '6.406'
$ws.Range("E51").Value = '  -7.65%  '
